$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Itens", "Itens - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 / column A currently holds the combined text; trim it down to
    # just the hosting-plan description.
    $ws.Range("A2").Value = "Starter Shared Hosting - (28/09/2016 - 27/10/2016)"

    # Push the old row 3 ("Late Fee" / "10.00") down to row 4, freeing up a
    # new row 3 for the "Hosting Location" detail line.
    $ws.Rows.Item(3).Insert()
    $ws.Range("A3").Value = "Hosting Location: United Kingdom (UK)"
}
